$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New header cells (bold style like existing headers A1:F1)
$ws.Range("G1").Value = "liberal"
$ws.Range("H1").Value = "kapica"
$ws.Range("I1").Value = "peker"

# Copy style from an existing header cell (F1) to the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New data cells in row 6
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 8

# Update selection to match diff (I7)
$ws.Range("I7").Select()
